# fix(publipostage): Try to solve Excel emoji problem
#
# Replace the three emoji "statut" markers with plain-text / warning
# equivalents everywhere they occur in the workbook:
#   📘 -> ⚠️
#   📙 -> +3
#   📕 -> -3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = [string]$cell.Value2

        if ($val -eq "📘") {
            $cell.Value = "⚠️"
        } elseif ($val -eq "📙") {
            # Force as literal text so Excel does not coerce "+3" into
            # the number 3 (dropping the leading plus sign).
            $origStyle = $cell.Style
            $cell.NumberFormat = "@"
            $cell.Value = "+3"
            $cell.Style = $origStyle
        } elseif ($val -eq "📕") {
            $origStyle = $cell.Style
            $cell.NumberFormat = "@"
            $cell.Value = "-3"
            $cell.Style = $origStyle
        }
    }
}
